$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 582023
$ws.Range("D2").Value = 153480
$ws.Range("E2").Value = 979264705
$ws.Range("C3").Value = 476
$ws.Range("E3").Value = 1001231
$ws.Range("C4").Value = 744
$ws.Range("E4").Value = 2622629
$ws.Range("C8").Value = 2603
$ws.Range("E8").Value = 14616027
$ws.Range("C10").Value = 242688
$ws.Range("D10").Value = 62741
$ws.Range("E10").Value = 982142970
$ws.Range("C11").Value = 650
$ws.Range("E11").Value = 14029671
$ws.Range("C12").Value = 37
$ws.Range("D12").Value = 15
$ws.Range("E12").Value = 1453120
$ws.Range("C13").Value = 128908
$ws.Range("D13").Value = 32007
$ws.Range("E13").Value = 550725129
$ws.Range("C16").Value = 7138
$ws.Range("D16").Value = 2786
$ws.Range("E16").Value = 15680255
$ws.Range("C18").Value = 200
$ws.Range("E18").Value = 1079889
$ws.Range("C19").Value = 16966
$ws.Range("D19").Value = 3902
$ws.Range("E19").Value = 56578997
$ws.Range("C21").Value = 135996
$ws.Range("D21").Value = 37457
$ws.Range("E21").Value = 228533543
$ws.Range("C27").Value = 64551
$ws.Range("E27").Value = 246816543
$ws.Range("C30").Value = 24513
$ws.Range("D30").Value = 6314
$ws.Range("E30").Value = 98851525
$ws.Range("C33").Value = 2981
$ws.Range("E33").Value = 8758324
$ws.Range("C35").Value = 4360
$ws.Range("D35").Value = 1016
$ws.Range("E35").Value = 13902065
$ws.Range("C36").Value = 164890
$ws.Range("D36").Value = 47550
$ws.Range("E36").Value = 289131493
$ws.Range("C39").Value = 3127
$ws.Range("E39").Value = 17965305
$ws.Range("C41").Value = 91011
$ws.Range("E41").Value = 371949138
$ws.Range("C44").Value = 20249
$ws.Range("D44").Value = 5305
$ws.Range("E44").Value = 96720089
$ws.Range("C46").Value = 2222
$ws.Range("E46").Value = 4734206
$ws.Range("C47").Value = 5663
$ws.Range("E47").Value = 19487330
$ws.Range("C48").Value = 118459
$ws.Range("D48").Value = 33071
$ws.Range("E48").Value = 200882904
$ws.Range("C52").Value = 1084
$ws.Range("E52").Value = 5291242
$ws.Range("C54").Value = 54774
$ws.Range("D54").Value = 14629
$ws.Range("E54").Value = 203373650
$ws.Range("C57").Value = 22522
$ws.Range("D57").Value = 5820
$ws.Range("E57").Value = 83987508
$ws.Range("C58").Value = 2435
$ws.Range("D58").Value = 849
$ws.Range("E58").Value = 6790269
$ws.Range("C60").Value = 3707
$ws.Range("D60").Value = 923
$ws.Range("E60").Value = 11093560
$ws.Range("C62").Value = 37037
$ws.Range("D62").Value = 9360
$ws.Range("E62").Value = 69566235
$ws.Range("C64").Value = 54
$ws.Range("E64").Value = 179346
$ws.Range("C66").Value = 18219
$ws.Range("E66").Value = 87570997
$ws.Range("C68").Value = 11919
$ws.Range("E68").Value = 52180439
$ws.Range("C69").Value = 973
$ws.Range("E69").Value = 2373317
$ws.Range("C70").Value = 691
$ws.Range("E70").Value = 1918721
$ws.Range("C71").Value = 254051
$ws.Range("D71").Value = 70478
$ws.Range("E71").Value = 446755467
$ws.Range("C72").Value = 277
$ws.Range("E72").Value = 1096467
$ws.Range("C75").Value = 1344
$ws.Range("E75").Value = 6831800
$ws.Range("C77").Value = 127197
$ws.Range("D77").Value = 34049
$ws.Range("E77").Value = 488730610
$ws.Range("C78").Value = 341
$ws.Range("E78").Value = 5685978
$ws.Range("C80").Value = 61368
$ws.Range("D80").Value = 16096
$ws.Range("E80").Value = 246271307
$ws.Range("C82").Value = 267
$ws.Range("D82").Value = 139
$ws.Range("E82").Value = 2475511
$ws.Range("C83").Value = 15193
$ws.Range("D83").Value = 6284
$ws.Range("E83").Value = 74643446
$ws.Range("C85").Value = 6911
$ws.Range("E85").Value = 22865174
$ws.Range("C86").Value = 51048
$ws.Range("D86").Value = 12032
$ws.Range("E86").Value = 78947816
$ws.Range("C87").Value = 16
$ws.Range("E87").Value = 22888
$ws.Range("C88").Value = 58
$ws.Range("E88").Value = 172627
$ws.Range("C89").Value = 12045
$ws.Range("D89").Value = 3175
$ws.Range("E89").Value = 25692041
$ws.Range("C91").Value = 11295
$ws.Range("E91").Value = 21757811
$ws.Range("C92").Value = 742
$ws.Range("E92").Value = 1241911
$ws.Range("C93").Value = 717
$ws.Range("E93").Value = 1227912
$ws.Range("C94").Value = 21064
$ws.Range("E94").Value = 42958261
$ws.Range("C97").Value = 6754
$ws.Range("E97").Value = 15942627
$ws.Range("C101").Value = 247210
$ws.Range("D101").Value = 66978
$ws.Range("E101").Value = 406134174
$ws.Range("C106").Value = 2929
$ws.Range("E106").Value = 14754756
$ws.Range("C108").Value = 103082
$ws.Range("D108").Value = 27268
$ws.Range("E108").Value = 385354170
$ws.Range("C109").Value = 332
$ws.Range("E109").Value = 5373964
$ws.Range("C111").Value = 56585
$ws.Range("D111").Value = 13939
$ws.Range("E111").Value = 216045827
$ws.Range("C114").Value = 2691
$ws.Range("D114").Value = 1034
$ws.Range("E114").Value = 7351019
$ws.Range("C115").Value = 5543
$ws.Range("D115").Value = 1311
$ws.Range("E115").Value = 16841275
$ws.Range("C117").Value = 994490
$ws.Range("D117").Value = 216796
$ws.Range("E117").Value = 1696145734
$ws.Range("C122").Value = 5102
$ws.Range("E122").Value = 44389947
$ws.Range("C124").Value = 432288
$ws.Range("D124").Value = 102266
$ws.Range("E124").Value = 1661257055
$ws.Range("C125").Value = 2010
$ws.Range("D125").Value = 544
$ws.Range("E125").Value = 35442676
$ws.Range("C127").Value = 398759
$ws.Range("D127").Value = 87171
$ws.Range("E127").Value = 1495645834
$ws.Range("C129").Value = 5036
$ws.Range("E129").Value = 9731610
$ws.Range("C131").Value = 16259
$ws.Range("D131").Value = 3937
$ws.Range("E131").Value = 52702187
$ws.Range("C134").Value = 61285
$ws.Range("D134").Value = 17282
$ws.Range("E134").Value = 88538978
$ws.Range("C139").Value = 18052
$ws.Range("E139").Value = 37486675
$ws.Range("C141").Value = 5143
$ws.Range("E141").Value = 10571237
$ws.Range("C144").Value = 493
$ws.Range("E144").Value = 792487
$ws.Range("C145").Value = 670
$ws.Range("E145").Value = 1276455
$ws.Range("C146").Value = 28256
$ws.Range("E146").Value = 43453023
$ws.Range("C149").Value = 11692
$ws.Range("D149").Value = 3136
$ws.Range("E149").Value = 29194752
$ws.Range("C151").Value = 8323
$ws.Range("E151").Value = 18567248
$ws.Range("C153").Value = 501
$ws.Range("E153").Value = 1071164
$ws.Range("C154").Value = 38428
$ws.Range("E154").Value = 93131657
$ws.Range("C156").Value = 486
$ws.Range("E156").Value = 1200074
$ws.Range("C159").Value = 152064
$ws.Range("D159").Value = 42023
$ws.Range("E159").Value = 264314332
$ws.Range("C164").Value = 2234
$ws.Range("E164").Value = 13088259
$ws.Range("C166").Value = 68482
$ws.Range("D166").Value = 18002
$ws.Range("E166").Value = 272729979
$ws.Range("C168").Value = 27930
$ws.Range("E168").Value = 118806346
$ws.Range("C170").Value = 2017
$ws.Range("E170").Value = 4335716
$ws.Range("C171").Value = 4299
$ws.Range("D171").Value = 1084
$ws.Range("E171").Value = 13840564
$ws.Range("C172").Value = 403179
$ws.Range("D172").Value = 113597
$ws.Range("E172").Value = 648853340
$ws.Range("C178").Value = 2374
$ws.Range("E178").Value = 12272168
$ws.Range("C180").Value = 169733
$ws.Range("D180").Value = 45422
$ws.Range("E180").Value = 654132829
$ws.Range("C181").Value = 381
$ws.Range("E181").Value = 6849705
$ws.Range("C183").Value = 69541
$ws.Range("D183").Value = 17645
$ws.Range("E183").Value = 280996685
$ws.Range("C186").Value = 9860
$ws.Range("D186").Value = 3681
$ws.Range("E186").Value = 31947787
$ws.Range("C188").Value = 11580
$ws.Range("D188").Value = 2783
$ws.Range("E188").Value = 33320402
$ws.Range("C190").Value = 462921
$ws.Range("D190").Value = 125085
$ws.Range("E190").Value = 721694802
$ws.Range("C198").Value = 191106
$ws.Range("D198").Value = 48841
$ws.Range("E198").Value = 714008512
$ws.Range("C199").Value = 404
$ws.Range("E199").Value = 7755574
$ws.Range("C201").Value = 108438
$ws.Range("D201").Value = 26284
$ws.Range("E201").Value = 401267829
$ws.Range("C204").Value = 8242
$ws.Range("D204").Value = 2935
$ws.Range("E204").Value = 19787719
$ws.Range("C207").Value = 14730
$ws.Range("D207").Value = 3325
$ws.Range("E207").Value = 41621345
$ws.Range("C209").Value = 180267
$ws.Range("D209").Value = 51681
$ws.Range("E209").Value = 299378902
$ws.Range("C213").Value = 1922
$ws.Range("E213").Value = 10736652
$ws.Range("C215").Value = 97380
$ws.Range("E215").Value = 378276912
$ws.Range("C218").Value = 26980
$ws.Range("D218").Value = 6928
$ws.Range("E218").Value = 119281726
$ws.Range("C220").Value = 3167
$ws.Range("E220").Value = 8527970
$ws.Range("C222").Value = 5917
$ws.Range("D222").Value = 1430
$ws.Range("E222").Value = 19237745
$ws.Range("C223").Value = 472071
$ws.Range("D223").Value = 119381
$ws.Range("E223").Value = 759456529
$ws.Range("C225").Value = 507
$ws.Range("E225").Value = 1679306
$ws.Range("C229").Value = 2712
$ws.Range("E229").Value = 16768462
$ws.Range("C231").Value = 201856
$ws.Range("D231").Value = 48839
$ws.Range("E231").Value = 791924534
$ws.Range("C232").Value = 396
$ws.Range("D232").Value = 136
$ws.Range("E232").Value = 8471783
$ws.Range("C234").Value = 144734
$ws.Range("D234").Value = 32979
$ws.Range("E234").Value = 554207412
$ws.Range("C237").Value = 4850
$ws.Range("D237").Value = 1639
$ws.Range("E237").Value = 13233407
$ws.Range("C240").Value = 11908
$ws.Range("D240").Value = 2578
$ws.Range("E240").Value = 35258858
